# Apply cryptos list price/volume refresh (GitHub Actions scheduled update)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '''59.219.60'
$ws.Range('E2').Value = '  +1.38%  '
$ws.Range('D3').Value = '''2.676.67'
$ws.Range('E3').Value = '  +5.61%  '
$ws.Range('E4').Value = '  +0.38%  '
$ws.Range('D5').Value = '''519.14'
$ws.Range('E5').Value = '  +2.58%  '
$ws.Range('D6').Value = '''145.61'
$ws.Range('E6').Value = '  +1.42%  '
$ws.Range('D7').Value = '''0.998'
$ws.Range('E7').Value = '  -0.06%  '
$ws.Range('D8').Value = '''0.571'
$ws.Range('E8').Value = '  +1.59%  '
$ws.Range('D9').Value = '''2.710.23'
$ws.Range('E9').Value = '  +6.87%  '
$ws.Range('D10').Value = '''6.28'
$ws.Range('E10').Value = '  +1.73%  '
$ws.Range('E11').Value = '  +4.70%  '
$ws.Range('E12').Value = '  +2.52%  '
$ws.Range('E13').Value = '  -1.82%  '
$ws.Range('D14').Value = '''3.152.08'
$ws.Range('E14').Value = '  +5.86%  '
$ws.Range('D15').Value = '''59.190.68'
$ws.Range('E15').Value = '  +1.39%  '
$ws.Range('D16').Value = '''21.14'
$ws.Range('E16').Value = '  +2.43%  '
$ws.Range('D17').Value = '''0.0000138'
$ws.Range('E17').Value = '  +2.37%  '
$ws.Range('D18').Value = '''2.712.12'
$ws.Range('E18').Value = '  +6.90%  '
$ws.Range('D19').Value = '''356.55'
$ws.Range('E19').Value = '  +6.47%  '
$ws.Range('D20').Value = '''4.57'
$ws.Range('E20').Value = '  +0.91%  '
$ws.Range('E21').Value = '  +4.06%  '
$ws.Range('E22').Value = '  +4.98%  '
$ws.Range('D23').Value = '''0.997'
$ws.Range('E23').Value = '  -0.16%  '
$ws.Range('D24').Value = '''62.12'
$ws.Range('E24').Value = '  +3.44%  '
$ws.Range('E25').Value = '  +3.89%  '
$ws.Range('B26').Value = 'Binance-PegBSC-USD'
$ws.Range('C26').Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range('D26').Value = '''0.993'
$ws.Range('E26').Value = '  -0.53%  '
$ws.Range('B27').Value = 'Kaspa'
$ws.Range('C27').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D27').Value = '''0.162'
$ws.Range('E27').Value = '  +1.24%  '
$ws.Range('D28').Value = '0.0₃0816'
$ws.Range('E28').Value = '  +3.85%  '
$ws.Range('D29').Value = '''7.25'
$ws.Range('E29').Value = '  +4.60%  '
$ws.Range('D30').Value = '''0.999'
$ws.Range('E30').Value = '  +0.00%  '
$ws.Range('D31').Value = '''6.38'
$ws.Range('E31').Value = '  +9.13%  '
$ws.Range('D32').Value = '''19.13'
$ws.Range('E32').Value = '  +3.48%  '
$ws.Range('E33').Value = '  +4.22%  '
$ws.Range('D34').Value = '''150.70'
$ws.Range('E34').Value = '  +0.57%  '
$ws.Range('D35').Value = '''0.980'
$ws.Range('E35').Value = '  +4.44%  '
$ws.Range('E36').Value = '  +3.46%  '
$ws.Range('E37').Value = '  +3.64%  '
$ws.Range('D38').Value = '''36.80'
$ws.Range('E38').Value = '  +2.10%  '
$ws.Range('E39').Value = '  +3.28%  '
$ws.Range('D40').Value = '''3.75'
$ws.Range('E40').Value = '  +6.44%  '
$ws.Range('E41').Value = '  +1.44%  '
$ws.Range('D42').Value = '''283.72'
$ws.Range('E42').Value = '  +0.14%  '
$ws.Range('D43').Value = '''0.620'
$ws.Range('E43').Value = '  +3.32%  '
$ws.Range('D44').Value = '''0.0993'
$ws.Range('E44').Value = '  -0.39%  '
$ws.Range('D45').Value = '''19.90'
$ws.Range('E45').Value = '  +6.79%  '
$ws.Range('D46').Value = '''0.992'
$ws.Range('E46').Value = '  -0.54%  '
$ws.Range('D47').Value = '''0.0534'
$ws.Range('E47').Value = '  +0.16%  '
$ws.Range('D48').Value = '''0.0233'
$ws.Range('E48').Value = '  +2.66%  '
$ws.Range('B49').Value = 'Maker'
$ws.Range('C49').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D49').Value = '''2.017.90'
$ws.Range('E49').Value = '  +7.01%  '
$ws.Range('B50').Value = 'RenderToken'
$ws.Range('C50').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D50').Value = '''4.74'
$ws.Range('E50').Value = '  +4.84%  '
$ws.Range('D51').Value = '''10.28'
$ws.Range('E51').Value = '  -0.04%  '
